# Daily attendance processing - 2026-01-17 19:32:21
# The "Recorded By" column (G) lists contributors; flip the ordering so the
# user's email is listed before "System" wherever both recorded a session.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.UsedRange
[void]$rng.Replace("System, dnasr281@gmail.com", "dnasr281@gmail.com, System")
